$wb = $excel.ActiveWorkbook

# --- Overall sheet -------------------------------------------------------
$overall = $wb.Worksheets.Item("Overall")

# B3 changes from a literal value to a formula (3.704 + 1.332); everything
# downstream (F3, G3, F15, G15) recalculates automatically.
$overall.Range("B3").Formula = "=3.704+1.332"

# Move the cursor/selection to E3 on the Overall sheet (was I19).
$overall.Activate()
$overall.Range("E3").Select()

# --- Electronics sheet -----------------------------------------------------
$electronics = $wb.Worksheets.Item("Electronics")

# Re-enter G2:G8's formula across the whole contiguous range so the engine
# recognizes it as a shared formula group (it previously was a set of
# individually-entered, non-shared formulas).
$electronics.Range("G2:G8").Formula = "=(E2/D2)*F2"

Write-Host "done"
